$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K" = strikeouts) values for rows 2-13 with newly
# regenerated figures (replacing the old Strike# based values).
$newK = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 2
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
